$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The practice list originally had 10 face/name/profession entries (rows 2-11).
# Two entries are removed entirely (their rows deleted, shifting everything
# below up): row 5 ("Linus" / Images_Faces/P4.jpg) and row 6 ("Tobias" /
# Images_Faces/P5.jpg).
$ws.Rows("5:6").Delete()

# After that delete, the former rows 10-11 ("Kemist"/P9.jpg and
# "Prast"/P10.jpg) are now rows 8-9. Remove those two as well, shrinking the
# fixation-cross / practice-trial list down to 6 entries total.
$ws.Rows("8:9").Delete()

# Restore the active selection to where the user left off editing.
$ws.Range("C13").Select()
